$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: description text now points at the local checkout path instead of
# the W:\ network share (fix_cbs_data_230717 folder).
$ws.Range("B2").Value = "C:\Users\dpere\Documents\JTMT\forecast\create_forecast_basic\current"

# B7: a spacer/blank row below the existing table (kept formatted so the
# row is materialised in the sheet).
$ws.Range("B7").NumberFormat = "General"

# B10: new cell holding (and linking to) the git-tools folder that used
# to live in B2, now exposed as a clickable hyperlink.
$ws.Range("B10").Value = "W:/Data/Forecast/Tools/forecast_git/create_forecast_basic/current"
$ws.Hyperlinks.Add($ws.Range("B10"), "W:/Data/Forecast/Tools/forecast_git/create_forecast_basic/current", "", "", "W:/Data/Forecast/Tools/forecast_git/create_forecast_basic/current")

# Hyperlinks.Add stamps the built-in blue/underlined "Hyperlink" style;
# put the plain black, non-underlined look back (matches the rest of the sheet).
$f = $ws.Range("B10").Font
$f.Underline = $false
$f.Color = 0

Write-Host "done"
